$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GAM Outputs")

# --- Update data rows 17-21 (column order: A Date, B Species, C Life History,
# D Model, E Deviance, F REML, G AIC, H s(lon,lat) edf, I s(doy) edf,
# P s(salinity) edf, Q s(temperature) edf, R s(salinity,temperature) edf) ---

# Row 17
$ws.Range("A17").Value = 44218
$ws.Range("E17").Value = 0.78800000000000003
$ws.Range("F17").Value = 2731.6
$ws.Range("G17").Value = "enter"
$ws.Range("H17").Value = 27.76
$ws.Range("I17").Value = 3.323

# Row 18
$ws.Range("A18").Value = 44218
$ws.Range("E18").Value = 0.78900000000000003
$ws.Range("F18").Value = 2732.9
$ws.Range("G18").Value = "enter"
$ws.Range("H18").Value = 27.760999999999999
$ws.Range("I18").Value = 3.2229999999999999
$ws.Range("P18").Value = 3.1589999999999998

# Row 19
$ws.Range("A19").Value = 44218
$ws.Range("E19").Value = 0.80600000000000005
$ws.Range("F19").Value = 2673.4
$ws.Range("G19").Value = "enter"
$ws.Range("H19").Value = 27.472000000000001
$ws.Range("I19").Value = 1.4710000000000001
$ws.Range("Q19").Value = 6.6879999999999997

# Row 20
$ws.Range("A20").Value = 44218
$ws.Range("E20").Value = 0.80900000000000005
$ws.Range("F20").Value = 2668.3
$ws.Range("G20").Value = "enter"
$ws.Range("H20").Value = 27.318999999999999
$ws.Range("I20").Value = 1
$ws.Range("P20").Value = 4.6040000000000001
$ws.Range("Q20").Value = 6.8239999999999998

# Row 21
$ws.Range("A21").Value = 44218
$ws.Range("E21").Value = 0.82699999999999996
$ws.Range("F21").Value = 2621.1
$ws.Range("G21").Value = "enter"
$ws.Range("H21").Value = 26.69
$ws.Range("I21").Value = 1.615
$ws.Range("R21").Value = 23.212

# --- Sheet view: scroll so column H is the leftmost visible column, and
# update the active selection to R22 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("R22").Select()
